$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.686.85'
$ws.Range("E2").Value = '  +0.38%  '
$ws.Range("D3").Value = '1.849.95'
$ws.Range("E3").Value = '  +0.10%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.08'
$ws.Range("E5").Value = '  +0.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9995'
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4270'
$ws.Range("E7").Value = '  +0.64%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3690'
$ws.Range("E8").Value = '  +1.43%  '
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07318'
$ws.Range("E9").Value = '  +0.42%  '
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8744'
$ws.Range("E10").Value = '  +0.21%  '
$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.94'
$ws.Range("E11").Value = '  +2.04%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.828.60'
$ws.Range("E12").Value = '  -0.19%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.435'
$ws.Range("E13").Value = '  +2.10%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.563'
$ws.Range("E14").Value = '  +0.68%  '
$ws.Range("B15").Value = 'TRON'
$ws.Range("C15").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06939'
$ws.Range("E15").Value = '  +0.56%  '
$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.002'
$ws.Range("E16").Value = '  -0.18%  '
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '80.65'
$ws.Range("E17").Value = '  +1.06%  '
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009075'
$ws.Range("E18").Value = '  +1.47%  '
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9984'
$ws.Range("E19").Value = '  -0.24%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.52'
$ws.Range("E20").Value = '  +1.12%  '
$ws.Range("B21").Value = 'WrappedBTC'
$ws.Range("C21").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D21").Value = '27.533.85'
$ws.Range("E21").Value = '  -0.25%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.086'
$ws.Range("E22").Value = '  +2.05%  '
$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.93'
$ws.Range("E23").Value = '  +5.60%  '
$ws.Range("B24").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D24").Value = '2.073.61'
$ws.Range("E24").Value = '  +1.50%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.959'
$ws.Range("E25").Value = '  -1.09%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.46'
$ws.Range("E26").Value = '  +0.17%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.51'
$ws.Range("E27").Value = '  -1.68%  '
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.269'
$ws.Range("E28").Value = '  +0.24%  '
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '115.28'
$ws.Range("E29").Value = '  -5.58%  '
$ws.Range("B30").Value = 'LidoDAOToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.860'
$ws.Range("E30").Value = '  -0.85%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08901'
$ws.Range("E31").Value = '  +0.29%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7811'
$ws.Range("E32").Value = '  +2.05%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.584'
$ws.Range("E33").Value = '  +1.04%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.972'
$ws.Range("E34").Value = '  -0.18%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.156'
$ws.Range("E35").Value = '  +4.73%  '
$ws.Range("B36").Value = 'Frax'
$ws.Range("C36").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9982'
$ws.Range("E36").Value = '  -0.25%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.106'
$ws.Range("E37").Value = '  +0.98%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05411'
$ws.Range("E38").Value = '  +0.79%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01950'
$ws.Range("E39").Value = '  +0.84%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.837'
$ws.Range("E40").Value = '  +0.43%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5150'
$ws.Range("E41").Value = '  +1.61%  '
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1669'
$ws.Range("E42").Value = '  +1.22%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.718'
$ws.Range("E43").Value = '  -2.47%  '
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.644'
$ws.Range("E44").Value = '  +3.63%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.65'
$ws.Range("E45").Value = '  +3.23%  '
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '107.33'
$ws.Range("E46").Value = '  +2.52%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.06561'
$ws.Range("E47").Value = '  +0.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4743'
$ws.Range("E48").Value = '  +1.40%  '
$ws.Range("B49").Value = 'PaxDollar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.9990'
$ws.Range("E49").Value = '  -0.20%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.655'
$ws.Range("E50").Value = '  +1.96%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.813'
$ws.Range("E51").Value = '  +5.43%  '
